# B1--and-B2-PowerPoint.pptx edit
#
# 1) Table on slide 5 (the "B1- TYPES OF FINANCIAL DOCUMENTS" slide) gets a
#    different built-in table style applied.
# 2) The deck's theme palette is swapped from the "Integral / Red Violet"
#    colours to the stock "Office" colour scheme.

$p = $ppt.ActivePresentation

# --- 1. Re-style the table on slide 5 -------------------------------------
$slide = $p.Slides.Item(5)
for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
    $shp = $slide.Shapes.Item($i)
    if ($shp.HasTable) {
        $shp.Table.ApplyStyle("{7FEA5404-F378-44B9-ADF0-A7A865311475}")
    }
}

# --- 2. Swap the theme colour scheme (Integral -> Office) -----------------
$design  = $p.Designs.Item(1)
$theme   = $design.SlideMaster.Theme
$colors  = $theme.ThemeColorScheme

# index -> (scheme slot, new RGB as a VBA-style BGR long == R + G*256 + B*65536)
$colors.Colors(1).RGB  = 0          # dk1      000000
$colors.Colors(2).RGB  = 16777215   # lt1      FFFFFF
$colors.Colors(3).RGB  = 6968388    # dk2      44546A
$colors.Colors(4).RGB  = 15132391   # lt2      E7E6E6
$colors.Colors(5).RGB  = 13998939   # accent1  5B9BD5
$colors.Colors(6).RGB  = 3243501    # accent2  ED7D31
$colors.Colors(7).RGB  = 10855845   # accent3  A5A5A5
$colors.Colors(8).RGB  = 49407      # accent4  FFC000
$colors.Colors(9).RGB  = 12874308   # accent5  4472C4
$colors.Colors(10).RGB = 4697456    # accent6  70AD47
$colors.Colors(11).RGB = 12673797   # hlink    0563C1
$colors.Colors(12).RGB = 7491477    # folHlink 954F72
